# Added DownLoad Issues File test
# Adds a new "testSuccessfulIssuesFileDownLoad" test-data block (header row +
# data row) to the TestCaseDataSets sheet, following the same layout used by
# the existing testSuccessfulLogIn / testSuccessfulLogOut blocks, but with two
# extra columns (downloadPath / filename). This fills the previously-blank
# rows 11-12, extends the blank rows below with the new D/E columns, and
# appends one more blank row (16) to keep the same "2 blank rows after each
# data block" spacing pattern used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCaseDataSets")

# --- Extend the still-blank rows 9 & 10 so they also carry blank D/E cells
#     (same style as their existing A:C cells) ------------------------------
$ws.Range("A9").Copy()
$ws.Range("D9:E9").PasteSpecial(-4122)

$ws.Range("A10").Copy()
$ws.Range("D10:E10").PasteSpecial(-4122)

# --- Row 11 / 12: new header + data block, styled like the
#     testSuccessfulLogOut block (rows 7-8) --------------------------------
$ws.Range("A7").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)

$ws.Range("A8").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)

$ws.Range("A11").Value = "testSuccessfulIssuesFileDownLoad"
$ws.Range("B11").Value = "username"
$ws.Range("C11").Value = "password"
$ws.Range("D11").Value = "downloadPath"
$ws.Range("E11").Value = "filename"

$ws.Range("A12").Value = "testSuccessfulIssuesFileDownLoad"
$ws.Range("B12").Value = "Richmond"
$ws.Range("C12").Value = "123456"
$ws.Range("E12").Value = "Richmond County.docx"
$ws.Range("D12").Value = "C:\\Users\\jslee\\Downloads"

# --- Rows 13-15 stay blank but now also carry blank D/E cells -------------
$ws.Range("A13").Copy()
$ws.Range("D13:E13").PasteSpecial(-4122)

$ws.Range("A14").Copy()
$ws.Range("D14:E14").PasteSpecial(-4122)

$ws.Range("A15").Copy()
$ws.Range("D15:E15").PasteSpecial(-4122)

# --- New trailing blank row 16, same style, spanning A:E -------------------
$ws.Range("A9").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column widths: widen D/E now that they hold the longer new strings ----
$ws.Columns.Item(4).ColumnWidth = 26.666666666666668
$ws.Columns.Item(5).ColumnWidth = 20.833333333333336

# --- Selection moves to D13, matching the post-edit cursor position --------
$ws.Range("D13").Select()
